# Added implementation of MSM measure.
# This adds the missing "inherited from java.lang.Object" operations (and the
# OrderController constructor) to the interfaceOperations sheet for the
# pl.piomin.order.controller.OrderController interface, and updates the two
# other sheets that referenced the "create()" operation signature for the row
# that now actually corresponds to the newly-added constructor signature.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. interfaceOperations sheet: expand pl.piomin.order.controller.OrderController
#    from its 4 declared operations to the full 13 operations (the inherited
#    java.lang.Object operations plus the constructor signature).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("interfaceOperations")

$ifaceName = "pl.piomin.order.controller.OrderController"

$ws2.Cells.Item(2, 1).Value = $ifaceName
$ws2.Cells.Item(2, 2).Value = "equals(java.lang.Object)"
$ws2.Cells.Item(2, 3).Value = "public"
$ws2.Cells.Item(2, 4).Value = "boolean"

$ws2.Cells.Item(3, 1).Value = $ifaceName
$ws2.Cells.Item(3, 2).Value = "create(pl.piomin.base.domain.Order)"
$ws2.Cells.Item(3, 3).Value = "public"
$ws2.Cells.Item(3, 4).Value = "pl.piomin.base.domain.Order"

$ws2.Cells.Item(4, 1).Value = $ifaceName
$ws2.Cells.Item(4, 2).Value = "toString()"
$ws2.Cells.Item(4, 3).Value = "public"
$ws2.Cells.Item(4, 4).Value = "java.lang.String"

$ws2.Cells.Item(5, 1).Value = $ifaceName
$ws2.Cells.Item(5, 2).Value = "all()"
$ws2.Cells.Item(5, 3).Value = "public"
$ws2.Cells.Item(5, 4).Value = "java.util.List"

$ws2.Cells.Item(6, 1).Value = $ifaceName
$ws2.Cells.Item(6, 2).Value = "OrderController(org.springframework.kafka.core.KafkaTemplate, org.springframework.kafka.config.StreamsBuilderFactoryBean, pl.piomin.order.service.OrderGeneratorService)"
$ws2.Cells.Item(6, 3).Value = "public"
$ws2.Cells.Item(6, 4).Value = "void"

$ws2.Cells.Item(7, 1).Value = $ifaceName
$ws2.Cells.Item(7, 2).Value = "getClass()"
$ws2.Cells.Item(7, 3).Value = "public"
$ws2.Cells.Item(7, 4).Value = "java.lang.Class"

$ws2.Cells.Item(8, 1).Value = $ifaceName
$ws2.Cells.Item(8, 2).Value = "notifyAll()"
$ws2.Cells.Item(8, 3).Value = "public"
$ws2.Cells.Item(8, 4).Value = "void"

$ws2.Cells.Item(9, 1).Value = $ifaceName
$ws2.Cells.Item(9, 2).Value = "hashCode()"
$ws2.Cells.Item(9, 3).Value = "public"
$ws2.Cells.Item(9, 4).Value = "int"

$ws2.Cells.Item(10, 1).Value = $ifaceName
$ws2.Cells.Item(10, 2).Value = "wait()"
$ws2.Cells.Item(10, 3).Value = "public"
$ws2.Cells.Item(10, 4).Value = "void"

$ws2.Cells.Item(11, 1).Value = $ifaceName
$ws2.Cells.Item(11, 2).Value = "notify()"
$ws2.Cells.Item(11, 3).Value = "public"
$ws2.Cells.Item(11, 4).Value = "void"

$ws2.Cells.Item(12, 1).Value = $ifaceName
$ws2.Cells.Item(12, 2).Value = "wait(long)"
$ws2.Cells.Item(12, 3).Value = "public"
$ws2.Cells.Item(12, 4).Value = "void"

$ws2.Cells.Item(13, 1).Value = $ifaceName
$ws2.Cells.Item(13, 2).Value = "create()"
$ws2.Cells.Item(13, 3).Value = "public"
$ws2.Cells.Item(13, 4).Value = "boolean"

$ws2.Cells.Item(14, 1).Value = $ifaceName
$ws2.Cells.Item(14, 2).Value = "wait(long, int)"
$ws2.Cells.Item(14, 3).Value = "public"
$ws2.Cells.Item(14, 4).Value = "void"

# ---------------------------------------------------------------------------
# 2. interfaceToClassRelations sheet: the relation row that used to point at
#    the "create()" signature actually originates from the OrderController
#    constructor signature.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("interfaceToClassRelations")
$ws6.Cells.Item(2, 2).Value = "OrderController(org.springframework.kafka.core.KafkaTemplate, org.springframework.kafka.config.StreamsBuilderFactoryBean, pl.piomin.order.service.OrderGeneratorService)"

# ---------------------------------------------------------------------------
# 3. methodNumberOfLines sheet: same correction - the 4-line method entry for
#    OrderController is really the constructor, not create().
# ---------------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item("methodNumberOfLines")
$ws11.Cells.Item(15, 2).Value = "OrderController(org.springframework.kafka.core.KafkaTemplate, org.springframework.kafka.config.StreamsBuilderFactoryBean, pl.piomin.order.service.OrderGeneratorService)"
